# Generate Report for Handback
# Update the generated/handoff/handback timestamps recorded on the
# Overview, zh-cn and de-de sheets to reflect a newer report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for 6b3cdc19-...md
$wsOverview.Range("G2").Value = "2016-10-19 17:25:38"

# zh-cn: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for 6b3cdc19-...md
$wsZhCn.Range("H2").Value = "2016-10-19 17:25:26"
$wsZhCn.Range("K2").Value = "2016-10-19 17:26:06"

# de-de: "Correspond Handoff Datetime" (matches Overview's generate date)
# and "Correspond Handback DateTime" for 6b3cdc19-...md
$wsDeDe.Range("H2").Value = "2016-10-19 17:25:38"
$wsDeDe.Range("K2").Value = "2016-10-19 17:26:24"
